# Updated cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.229.16'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.782.10'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.36'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.547'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '31.90'
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.039.62'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.02'
$ws.Range('E13').Value = '  -1.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.771.38'
$ws.Range('E14').Value = '  -0.65%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '34.157.46'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.623'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '67.90'
$ws.Range('E18').Value = '  +1.91%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '247.16'
$ws.Range('E19').Value = '  +3.78%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  +3.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.99'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '162.40'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  +2.45%  '
$ws.Range('E27').Value = '  +1.38%  '
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('E31').Value = '  +1.96%  '
$ws.Range('E32').Value = '  +4.35%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.70'
$ws.Range('E33').Value = '  +5.74%  '
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.444.73'
$ws.Range('E35').Value = '  +4.48%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.653'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.42'
$ws.Range('E37').Value = '  +7.95%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0192'
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.04'
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.39'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '80.22'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.923'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.64'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('E44').Value = '  +0.68%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.08'
$ws.Range('E45').Value = '  +3.96%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.940.85'
$ws.Range('E49').Value = '  +0.14%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '104.77'
